# Update "Liver vascular" term to "Liver vasculature" and add a new
# YouTube hyperlink for the "Portal vein thrombosis" row (row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the term "Liver vascular" -> "Liver vasculature" (affects A10 & A11
# since they share the same text).
$ws.Range("A10").Value = "Liver vasculature"
$ws.Range("A11").Value = "Liver vasculature"

# Add the new hyperlink in D10 pointing to the new YouTube video.
$ws.Hyperlinks.Add($ws.Range("D10"), "https://youtu.be/DjI1kEnzfSQ ") | Out-Null
$ws.Range("D10").Style = $ws.Range("D3").Style

# Update the current selection to D12 (matches the saved view state).
$ws.Range("D12").Select()
